$d = $word.ActiveDocument

# 1. "Larve ;" -> "Larves ;"
$d.Content.Find.Execute("Larve ;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Larves ;", 2)

# 2. "Pupae (besoin de trouver une traduction) ;" -> "Nymphes ;"
$d.Content.Find.Execute("Pupae (besoin de trouver une traduction) ;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Nymphes ;", 2)

# 3. "Adultes gesting (besoin aussi de trouver une trad) ;" -> "Adultes en gestation ;"
$d.Content.Find.Execute("Adultes gesting (besoin aussi de trouver une trad) ;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Adultes en gestation ;", 2)

# 4. "Adultes reproductifs." -> "Adultes capable de se reproduire."
$d.Content.Find.Execute("Adultes reproductifs.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Adultes capable de se reproduire.", 2)

# 5. Rewrite the long paragraph about serotypes, which is split across many runs.
#    Find the paragraph containing the unique marker text, delete its whole
#    range (minus the paragraph mark) and insert the replacement as a single run.
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*L’hypothèse la plus contraignante pour le modèle est la dernière*") {
        $r = $p.Range
        [void]$r.MoveEnd(1, -1)
        $r.Delete()
        $r.InsertAfter("L’hypothèse la plus contraignante pour le modèle est la dernière. Il existe en réalité plusieurs sérotypes pour ce virus. Hors, si une personne a déjà été infectée par le passé elle devrait se trouver dans le compartiment immunisés et non pas succeptible. De plus, une personne déjà infectée par un serotype peut être infectée par un autre serotype. Si cette hypothèse était relaxée, elle complefirai grandement les équations consernant les compartiments succeptibles et immunisés.")
        break
    }
}
